$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated / new data rows for the Vwf-Tnfrsf11b LR-pairs table (rows 2-10)
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric NATMI metrics recomputed with new TPM values.
$data = @(
    @("ECs", "Vwf", "Tnfrsf11b", "ECs", 3, 1, 18.975105, 56.925315, 0.9552145540969871, 0.955214554096987, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 1.438471084875, 12.946239763875, 0.02692841337005778, 0.02692841337005778),
    @("ECs", "Vwf", "Tnfrsf11b", "FAPs", 3, 1, 18.975105, 56.925315, 0.9552145540969871, 0.955214554096987, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 39.756557770545, 357.8090199349049, 0.7442492470461094, 0.7442492470461093),
    @("ECs", "Vwf", "Tnfrsf11b", "MuSCs", 3, 1, 18.975105, 56.925315, 0.9552145540969871, 0.955214554096987, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 9.830944975184998, 88.47850477666499, 0.18403689368082, 0.18403689368082),
    @("FAPs", "Vwf", "Tnfrsf11b", "ECs", 1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.01570916103663723, 0.01570916103663723, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 0.02365664742222222, 0.2129098268, 0.000442856298908969, 0.000442856298908969),
    @("FAPs", "Vwf", "Tnfrsf11b", "FAPs", 1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.01570916103663723, 0.01570916103663723, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 0.6538239661457778, 5.884415695312, 0.01223969130610236, 0.01223969130610236),
    @("FAPs", "Vwf", "Tnfrsf11b", "MuSCs", 1, 0.3333333333333333, 0.3120586666666667, 0.936176, 0.01570916103663723, 0.01570916103663723, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 0.161676659024, 1.455089931216, 0.0030266134316259, 0.0030266134316259),
    @("MuSCs", "Vwf", "Tnfrsf11b", "ECs", 3, 1, 0.5775933333333333, 1.73278, 0.02907628486637583, 0.02907628486637583, 1, 0.3333333333333333, 0.07580833333333332, 0.227425, 0.028190957994264, 0.02819095799426401, 0.04378638794444444, 0.3940774915, 0.0008196883252972552, 0.0008196883252972552),
    @("MuSCs", "Vwf", "Tnfrsf11b", "FAPs", 3, 1, 0.5775933333333333, 1.73278, 0.02907628486637583, 0.02907628486637583, 3, 1, 2.095195666666667, 6.285587, 0.7791435378093522, 0.7791435378093522, 1.210171049095556, 10.89153944186, 0.02265459945714059, 0.02265459945714059),
    @("MuSCs", "Vwf", "Tnfrsf11b", "MuSCs", 3, 1, 0.5775933333333333, 1.73278, 0.02907628486637583, 0.02907628486637583, 3, 1, 0.5180969999999999, 1.554291, 0.1926655041963838, 0.1926655041963838, 0.2992493732199999, 2.69324435898, 0.005601997083937984, 0.005601997083937984),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
